$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and add a new
# ListBullet paragraph right after it naming the responsible instructor.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $targetIndex = $target.Index
    $target.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "210064 - Eduardo Rezende Triboni"
    $newPara.Style = "ListBullet"
}

Write-Output "done"
